# Creating special classes to handle the different types of cards.
# Adds five new sheets (Special Cards, Transportation, Energy, Community Cards,
# Chance Cards) alongside the existing "European Cities" sheet, populates them
# with Monopoly reference data, and updates view/selection state to match.

$wb = $excel.ActiveWorkbook

function Set-RowValues($ws, [int]$row, [object[]]$values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $val = $values[$i]
        if ($null -ne $val) {
            $ws.Cells.Item($row, $i + 1).Value = $val
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet 1: European Cities (existing sheet) - selection moves, tab deselected
# ---------------------------------------------------------------------------
$wsCities = $wb.Worksheets.Item(1)
$wsCities.Range("D54").Select()

# ---------------------------------------------------------------------------
# Sheet 2: Special Cards
# ---------------------------------------------------------------------------
$wsSpecial = $wb.Worksheets.Add($null, $wsCities)
$wsSpecial.Name = "Special Cards"

Set-RowValues $wsSpecial 1 @("City", "Money Get", "Money Lost", "Action")
Set-RowValues $wsSpecial 2 @("Starting Point", 200000, $null, "Get Money")
Set-RowValues $wsSpecial 3 @("Visit Prison", 0, $null, "None")
Set-RowValues $wsSpecial 4 @("Free Holidays", 0, $null, "None")
Set-RowValues $wsSpecial 5 @("Go to Jail", 0, $null, "Go to jail")
Set-RowValues $wsSpecial 6 @("Income Tax", 0, 100000, "Pay")
Set-RowValues $wsSpecial 7 @("Super Tax", $null, 200000, "Pay")

$wsSpecial.Range("A1:F1").EntireColumn.ColumnWidth = 17.6640625
$wsSpecial.Range("D8").Select()

# ---------------------------------------------------------------------------
# Sheet 3: Transportation
# ---------------------------------------------------------------------------
$wsTransport = $wb.Worksheets.Add($null, $wsSpecial)
$wsTransport.Name = "Transportation"

Set-RowValues $wsTransport 1 @("City", "Sell Value", "Resell Value", "Rent")
Set-RowValues $wsTransport 2 @("Paris to Bourdeux", 250000, 75000, 40000)
Set-RowValues $wsTransport 3 @("Coppenhagen to Malmo", 250000, 75000, 40000)
Set-RowValues $wsTransport 4 @("London to Garforth", 250000, 75000, 40000)
Set-RowValues $wsTransport 5 @("Vienna to Walserberg", 250000, 75000, 40000)

$wsTransport.Range("A1").EntireColumn.ColumnWidth = 21.88671875
$wsTransport.Range("B1:F1").EntireColumn.ColumnWidth = 17.6640625
$wsTransport.Range("O21").Select()

# ---------------------------------------------------------------------------
# Sheet 4: Energy
# ---------------------------------------------------------------------------
$wsEnergy = $wb.Worksheets.Add($null, $wsTransport)
$wsEnergy.Name = "Energy"

Set-RowValues $wsEnergy 1 @("City", "Sell Value", "Resell Value", "Rent")
Set-RowValues $wsEnergy 2 @("Hydro Electric", 250000, 75000, 40000)
Set-RowValues $wsEnergy 3 @("Wind Electricity", 250000, 75000, 40000)

$wsEnergy.Range("A1").EntireColumn.ColumnWidth = 14.5546875
$wsEnergy.Range("B1").EntireColumn.ColumnWidth = 9.5546875
$wsEnergy.Range("C1").EntireColumn.ColumnWidth = 11.6640625
$wsEnergy.Range("D1").EntireColumn.ColumnWidth = 6
$wsEnergy.Range("H18").Select()

# ---------------------------------------------------------------------------
# Sheet 5: Community Cards
# ---------------------------------------------------------------------------
$wsCommunity = $wb.Worksheets.Add($null, $wsEnergy)
$wsCommunity.Name = "Community Cards"

Set-RowValues $wsCommunity 1 @("Card", "Message")
Set-RowValues $wsCommunity 2 @("Community", "Advance to the starting point")
Set-RowValues $wsCommunity 3 @("Community", "Bank Error: You get 200,000 from the Bank")
Set-RowValues $wsCommunity 4 @("Community", "Get out of jail free card")
Set-RowValues $wsCommunity 5 @("Community", "Go to prison without passig through the starting point")
Set-RowValues $wsCommunity 6 @("Community", "Happy birthday, get 20,000 from each player")
Set-RowValues $wsCommunity 7 @("Community", "Life Insurance: Collect 50,000 from the bank")
Set-RowValues $wsCommunity 8 @("Community", "School fees: Pay 50,000 for school tuition")
Set-RowValues $wsCommunity 9 @("Community", "You have won the 2nd place in an ugliness competition. Get 10,000")
Set-RowValues $wsCommunity 10 @("Community", "You inherit 100,000 from your aunt")
Set-RowValues $wsCommunity 11 @("Community", "Grand Opera! You need a new tux. Pay 20,000")

$wsCommunity.Range("A1").EntireColumn.ColumnWidth = 18.21875
$wsCommunity.Range("B1").EntireColumn.ColumnWidth = 59.44140625
$wsCommunity.Range("B10").Select()

# ---------------------------------------------------------------------------
# Sheet 6: Chance Cards
# ---------------------------------------------------------------------------
$wsChance = $wb.Worksheets.Add($null, $wsCommunity)
$wsChance.Name = "Chance Cards"

Set-RowValues $wsChance 1 @("Card", "Message")
Set-RowValues $wsChance 2 @("Chance", "Advance to the starting point")
Set-RowValues $wsChance 3 @("Chance", "Take a walk. Advance 5 blocks")
Set-RowValues $wsChance 4 @("Chance", "Bank offers you dividends. Collect 20,000")
Set-RowValues $wsChance 5 @("Chance", "Pay poor taxes of 25,000")
Set-RowValues $wsChance 6 @("Chance", "Get out of jail free card")
Set-RowValues $wsChance 7 @("Chance", "Go to jail without passing from the starting point")
Set-RowValues $wsChance 8 @("Chance", "Go back 3 boxes")
Set-RowValues $wsChance 9 @("Chance", "Your houses collapsed. Pay 25,000 for each house and hotel for repears")
Set-RowValues $wsChance 10 @("Chance", "Go for a reading. Advance 10 blocks")
Set-RowValues $wsChance 11 @("Chance", "You broke your leg sking. Pay 100,000 for medical expenses")

$wsChance.Range("B1").EntireColumn.ColumnWidth = 63.21875
$wsChance.Range("F25").Select()

# ---------------------------------------------------------------------------
# Final active sheet/tab: Transportation (matches bookViews activeTab="2")
# ---------------------------------------------------------------------------
$wsTransport.Activate()
$wsTransport.Range("O21").Select()
